$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mass")

# Update FOVT STATUS column values from "in fovt" to "in FOVT" for rows 2-5 and 7
$ws.Range("A2").Value = "in FOVT"
$ws.Range("A3").Value = "in FOVT"
$ws.Range("A4").Value = "in FOVT"
$ws.Range("A5").Value = "in FOVT"
$ws.Range("A7").Value = "in FOVT"

# Row 9: new "chest circumference" term, in OBA, synonym "girth"
$ws.Range("A9").Value = "in OBA"
$ws.Range("C9").Value = "chest circumference"
$ws.Range("D9").Value = "girth"

# Update the active selection to C7
$ws.Activate()
$ws.Range("C7").Select()
